$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the default (unstyled) cell style to restore after forcing text format,
# since assigning numeric-looking strings via .Value can cause Excel to coerce them
# into numbers and silently apply a new number-format style to the cell.
$origStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "27.518.58"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").Value = "1.618.57"
$ws.Range("E3").Value = "  -1.26%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.99"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -0.70%  "

$ws.Range("E6").Value = "  -1.85%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.76"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -1.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +1.86%  "

$ws.Range("E11").Value = "  -0.45%  "

$ws.Range("D12").Value = "1.847.89"
$ws.Range("E12").Value = "  -1.25%  "

$ws.Range("D13").Value = "1.622.87"
$ws.Range("E13").Value = "  -1.38%  "

$ws.Range("E14").Value = "  -0.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.550"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -2.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.58"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +0.84%  "

$ws.Range("D17").Value = "27.524.79"
$ws.Range("E17").Value = "  -0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.88"
$ws.Range("D18").Style = $origStyle

$ws.Range("E19").Value = "  -0.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.55"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -1.56%  "

$ws.Range("E22").Value = "  -0.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.96"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  +0.82%  "

$ws.Range("E24").Value = "  +6.55%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.90"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -0.26%  "

$ws.Range("E26").Value = "  -1.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.82"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -1.88%  "

$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.56"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("E30").Value = "  -0.47%  "

$ws.Range("E31").Value = "  -0.96%  "

$ws.Range("E32").Value = "  -0.38%  "

$ws.Range("D33").Value = "1.446.50"
$ws.Range("E33").Value = "  +1.41%  "

$ws.Range("E34").Value = "  -2.90%  "

$ws.Range("E35").Value = "  -2.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.944"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +4.88%  "

$ws.Range("E38").Value = "  -1.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0168"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +0.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.862"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -1.64%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.50"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  +7.05%  "

$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("E43").Value = "  -1.94%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("E45").Value = "  -2.07%  "

$ws.Range("E46").Value = "  -2.12%  "

$ws.Range("D47").Value = "1.759.19"
$ws.Range("E47").Value = "  -1.22%  "

$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.36"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0990"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +0.52%  "
